# Weekly update: insert a new price record as row 6, shifting the
# existing rows 6-12 down to 7-13 (dimension grows from A1:T12 to A1:T13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 6 - this shifts rows 6:12 down
# to 7:13 and keeps all their original data/styles intact.
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with the new weekly record. Columns A-L mirror
# the other rows for this market/product (same header dimensions).
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 44935
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100101
$ws.Range("H6").Value = "Berries"
$ws.Range("I6").Value = 100101001
$ws.Range("J6").Value = "Arándano (blue)"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 3000
$ws.Range("O6").Value = 3000
$ws.Range("P6").Value = 3000
$ws.Range("Q6").Value = "$/bandeja 2 kilos"
$ws.Range("R6").Value = "Provincia de Diguillín"
$ws.Range("S6").Value = 1500
$ws.Range("T6").Value = 2
